# Updates cryptocurrency price (column D) and 1h-volume-change (column E)
# figures on the "cryptos" worksheet to a newer snapshot, per the
# automated "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Price-looking numeric strings (e.g. "217.83")
# are prefixed with a leading apostrophe so Excel stores them as literal
# text (matching the original inline-string cells) instead of silently
# re-parsing them into numbers and dropping formatting such as trailing
# zeros (e.g. "0.07380" -> 0.0738).
$updates = [ordered]@{
    'D2' = '26.232.04'
    'E2' = '  -1.67%  '
    'D3' = '1.672.32'
    'E3' = '  +1.11%  '
    'D4' = '''1.007'
    'E4' = '  -0.19%  '
    'D5' = '''217.83'
    'E5' = '  -0.87%  '
    'D6' = '''0.5129'
    'E6' = '  +0.44%  '
    'D7' = '''1.007'
    'E7' = '  -0.28%  '
    'D8' = '''0.2658'
    'E8' = '  +4.83%  '
    'D9' = '''0.06385'
    'E9' = '  +4.11%  '
    'D10' = '''21.53'
    'E10' = '  -0.65%  '
    'D11' = '''0.07380'
    'E11' = '  +0.19%  '
    'D12' = '1.676.20'
    'E12' = '  +1.44%  '
    'D13' = '''4.552'
    'E13' = '  +2.12%  '
    'D14' = '''0.5828'
    'E14' = '  +1.76%  '
    'D15' = '1.902.12'
    'E15' = '  +1.46%  '
    'D16' = '''0.000008686'
    'E16' = '  +7.50%  '
    'D17' = '''64.76'
    'E17' = '  +0.38%  '
    'D18' = '26.301.46'
    'E18' = '  -1.24%  '
    'D19' = '''4.954'
    'E19' = '  -0.40%  '
    'E20' = '  -0.51%  '
    'D21' = '''10.85'
    'E21' = '  +3.06%  '
    'D22' = '''189.38'
    'E22' = '  +4.11%  '
    'D23' = '''6.213'
    'E23' = '  +0.40%  '
    'E24' = '  -0.26%  '
    'D25' = '''144.60'
    'E25' = '  +1.01%  '
    'D26' = '''7.627'
    'E26' = '  +0.55%  '
    'D27' = '''0.1184'
    'E27' = '  +3.98%  '
    'D28' = '''15.64'
    'E28' = '  +4.23%  '
    'D29' = '''0.05928'
    'E29' = '  +2.45%  '
    'D30' = '''1.282'
    'E30' = '  -3.50%  '
    'E31' = '  -1.09%  '
    'D32' = '''3.525'
    'E32' = '  +2.83%  '
    'D33' = '''3.527'
    'E33' = '  +3.41%  '
    'D34' = '''1.637'
    'E34' = '  +3.80%  '
    'D35' = '''1.014'
    'E35' = '  +3.15%  '
    'D36' = '''0.6030'
    'E36' = '  +1.05%  '
    'D37' = '''2.373'
    'E37' = '  -2.33%  '
    'D38' = '''2.649'
    'E38' = '  +0.65%  '
    'D39' = '''0.01618'
    'E39' = '  +3.06%  '
    'D40' = '''6.078'
    'E40' = '  +6.40%  '
    'D41' = '1.078.59'
    'E41' = '  +0.75%  '
    'D42' = '''0.8715'
    'E42' = '  +0.83%  '
    'D43' = '''1.011'
    'E43' = '  -0.09%  '
    'D44' = '''99.95'
    'E44' = '  +4.42%  '
    'D45' = '1.822.54'
    'E45' = '  +1.89%  '
    'D46' = '''0.00000000114'
    'E46' = '  +5.80%  '
    'D47' = '''56.06'
    'E48' = '  -0.55%  '
    'D49' = '''8.088'
    'E49' = '  +4.17%  '
    'D50' = '''0.4304'
    'E50' = '  -1.82%  '
    'E51' = '  +0.01%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
